$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A2:G59")
$key1 = $ws.Range("D2:D59")
$key2 = $ws.Range("A2:A59")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1, 0, 1, 0, 0) | Out-Null
$ws.Sort.SortFields.Add($key2, 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

$ws.Range("C60").Select()
